# Append two new job postings (2026-01-21 12:44 JST) into the
# "ランサーズ" sheet of 案件情報.xlsx:
#   - row 5: "継続依頼あり" education-system closing/requirements job
#   - row 7 (after the shift caused by the first insert): auto-repair CRM job
# All existing rows keep their data but get the refreshed fetch timestamp,
# and the rows below each insertion point shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert the two new rows, shifting everything below down ---
# First insert creates a blank row 5 (old row5 "野球..." -> row6, old row6 "CSV..." -> row7)
$ws.Rows.Item(5).Insert()
# Second insert creates a blank row 7 (old row6 "CSV..." now at row7 -> row8)
$ws.Rows.Item(7).Insert()

# --- Refresh the "取得日時" timestamp for every data row ---
for ($r = 2; $r -le 8; $r++) {
  $ws.Cells.Item($r, 1).Value = "2026-01-21 12:44:19"
}

# --- New row 5: continuing-education system closing / requirements job ---
$ws.Cells.Item(5, 2).Value = "【継続依頼あり】教育システム開発案件のクロージング代行&要件定義"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5475937"
$ws.Cells.Item(5, 7).Value = 118

# --- New row 7: auto-repair shop CRM build partner job ---
$ws.Cells.Item(7, 2).Value = "【急募】自動車整備業向けCRM構築パートナー募集"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5474125"
$ws.Cells.Item(7, 7).Value = 25

# --- Column D (価格) gets a bit wider ---
# NB: the ColumnWidth COM property adds ~0.8333 of padding relative to the
# raw OOXML <col width> units, so back that padding out to land on exactly 30.
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668

# --- Rebuild the hyperlinks on column F top to bottom so every F cell
#     (old and new) points at the right URL with no stale/duplicate links ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5455098")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5445159")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5445154")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5475937")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5475665")
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://www.lancers.jp/work/detail/5474125")
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), "https://www.lancers.jp/work/detail/5475924")

# Hyperlinks.Add() re-styles the cell with a fresh (duplicate) style; put the
# normal "Hyperlink" cell style back on every URL cell, matching the original.
for ($r = 2; $r -le 8; $r++) {
  $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
